# edit.ps1 — applies the "Tolerance, barcode search, report update" changes
# to the 14575A00 Raw Data screening report table.
#
# 1) Chemistry line: "Zn/MnO2" -> "Zn/MnO₂" (subscript 2, U+2082).
# 2) Raw-data table: per-cell updates (Inspection column Pass -> OK, except
#    the two rows that become "Tolerance"; the Outlier row becomes
#    "OutlierL"; and specific OCV/CCV numeric cells get a trailing
#    " !" / " ^" / " *" annotation) as captured from the source diff.

$d = $word.ActiveDocument

# --- 1. Fix the chemistry label's subscript ---------------------------------
$subscript2 = [char]0x2082
$null = $d.Content.Find.Execute(
    "Zn/MnO2", $false, $false, $false, $false, $false, $true, 1, $false,
    "Zn/MnO" + $subscript2, 2)

# --- 2. Raw-data table cell updates ------------------------------------------
$t = $d.Tables.Item(1)

$cellUpdates = @(
    @{Row=2; Col=6; Text="1.572 ^"},
    @{Row=2; Col=8; Text="Tolerance"},
    @{Row=3; Col=8; Text="OK"},
    @{Row=4; Col=8; Text="OK"},
    @{Row=5; Col=8; Text="OK"},
    @{Row=6; Col=8; Text="OK"},
    @{Row=7; Col=8; Text="OK"},
    @{Row=8; Col=6; Text="1.568 ^"},
    @{Row=8; Col=8; Text="Tolerance"},
    @{Row=9; Col=8; Text="OK"},
    @{Row=10; Col=8; Text="OK"},
    @{Row=11; Col=8; Text="OK"},
    @{Row=12; Col=8; Text="OK"},
    @{Row=13; Col=8; Text="OK"},
    @{Row=14; Col=8; Text="OK"},
    @{Row=15; Col=8; Text="OK"},
    @{Row=16; Col=8; Text="OK"},
    @{Row=17; Col=4; Text="1.565 !"},
    @{Row=17; Col=6; Text="1.565 !"},
    @{Row=18; Col=8; Text="OK"},
    @{Row=19; Col=8; Text="OK"},
    @{Row=20; Col=8; Text="OK"},
    @{Row=21; Col=8; Text="OK"},
    @{Row=22; Col=8; Text="OK"},
    @{Row=23; Col=8; Text="OK"},
    @{Row=24; Col=8; Text="OK"},
    @{Row=25; Col=8; Text="OK"},
    @{Row=26; Col=4; Text="1.569 !"},
    @{Row=26; Col=6; Text="1.569 !"},
    @{Row=26; Col=7; Text="1.519 *"},
    @{Row=27; Col=8; Text="OK"},
    @{Row=28; Col=4; Text="1.569 !"},
    @{Row=29; Col=8; Text="OK"},
    @{Row=30; Col=4; Text="1.554 !"},
    @{Row=30; Col=6; Text="1.554 !"},
    @{Row=30; Col=7; Text="1.552 *"},
    @{Row=31; Col=8; Text="OK"},
    @{Row=32; Col=4; Text="1.565 !"},
    @{Row=32; Col=6; Text="1.566 !"},
    @{Row=32; Col=7; Text="1.518 *"},
    @{Row=33; Col=8; Text="OK"},
    @{Row=34; Col=8; Text="OK"},
    @{Row=35; Col=8; Text="OK"},
    @{Row=36; Col=6; Text="1.569 !"},
    @{Row=36; Col=7; Text="1.567 *"},
    @{Row=37; Col=8; Text="OK"},
    @{Row=38; Col=8; Text="OK"},
    @{Row=39; Col=8; Text="OK"},
    @{Row=40; Col=7; Text="1.497 *"},
    @{Row=40; Col=8; Text="OutlierL"},
    @{Row=41; Col=8; Text="OK"},
    @{Row=42; Col=4; Text="1.561 !"},
    @{Row=42; Col=6; Text="1.561 !"},
    @{Row=42; Col=7; Text="1.514 *"},
    @{Row=43; Col=8; Text="OK"},
    @{Row=44; Col=4; Text="1.567 !"},
    @{Row=44; Col=6; Text="1.568 !"},
    @{Row=45; Col=8; Text="OK"},
    @{Row=46; Col=8; Text="OK"},
    @{Row=47; Col=8; Text="OK"},
    @{Row=48; Col=8; Text="OK"},
    @{Row=49; Col=8; Text="OK"},
    @{Row=50; Col=8; Text="OK"},
    @{Row=51; Col=8; Text="OK"}
)

foreach ($u in $cellUpdates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $cell.Range.Text = $u.Text
}

Write-Host ("Updated " + $cellUpdates.Count + " table cells.")
